$d = $word.ActiveDocument

$d.Content.Find.Execute("638÷9=70, 8", $true, $false, $false, $false, $false, $true, 1, $false, "191÷3=63, 2", 2) | Out-Null
$d.Content.Find.Execute("282÷6=47, 0", $true, $false, $false, $false, $false, $true, 1, $false, "901÷8=112, 5", 2) | Out-Null
$d.Content.Find.Execute("615÷5=123, 0", $true, $false, $false, $false, $false, $true, 1, $false, "319÷8=39, 7", 2) | Out-Null
$d.Content.Find.Execute("123÷4=30, 3", $true, $false, $false, $false, $false, $true, 1, $false, "432÷8=54, 0", 2) | Out-Null
$d.Content.Find.Execute("537÷8=67, 1", $true, $false, $false, $false, $false, $true, 1, $false, "841÷5=168, 1", 2) | Out-Null
$d.Content.Find.Execute("765÷6=127, 3", $true, $false, $false, $false, $false, $true, 1, $false, "771÷2=385, 1", 2) | Out-Null
$d.Content.Find.Execute("723÷4=180, 3", $true, $false, $false, $false, $false, $true, 1, $false, "771÷3=257, 0", 2) | Out-Null
$d.Content.Find.Execute("821÷5=164, 1", $true, $false, $false, $false, $false, $true, 1, $false, "800÷9=88, 8", 2) | Out-Null
$d.Content.Find.Execute("150÷5=30, 0", $true, $false, $false, $false, $false, $true, 1, $false, "932÷5=186, 2", 2) | Out-Null
$d.Content.Find.Execute("713÷3=237, 2", $true, $false, $false, $false, $false, $true, 1, $false, "639÷9=71, 0", 2) | Out-Null
$d.Content.Find.Execute("564÷2=282, 0", $true, $false, $false, $false, $false, $true, 1, $false, "640÷5=128, 0", 2) | Out-Null
$d.Content.Find.Execute("943÷9=104, 7", $true, $false, $false, $false, $false, $true, 1, $false, "614÷7=87, 5", 2) | Out-Null
$d.Content.Find.Execute("486÷4=121, 2", $true, $false, $false, $false, $false, $true, 1, $false, "702÷7=100, 2", 2) | Out-Null
$d.Content.Find.Execute("710÷9=78, 8", $true, $false, $false, $false, $false, $true, 1, $false, "394÷9=43, 7", 2) | Out-Null
$d.Content.Find.Execute("795÷3=265, 0", $true, $false, $false, $false, $false, $true, 1, $false, "112÷2=56, 0", 2) | Out-Null
$d.Content.Find.Execute("689÷2=344, 1", $true, $false, $false, $false, $false, $true, 1, $false, "758÷3=252, 2", 2) | Out-Null
$d.Content.Find.Execute("791÷2=395, 1", $true, $false, $false, $false, $false, $true, 1, $false, "876÷9=97, 3", 2) | Out-Null
$d.Content.Find.Execute("632÷7=90, 2", $true, $false, $false, $false, $false, $true, 1, $false, "449÷6=74, 5", 2) | Out-Null
$d.Content.Find.Execute("897÷4=224, 1", $true, $false, $false, $false, $false, $true, 1, $false, "340÷7=48, 4", 2) | Out-Null
$d.Content.Find.Execute("839÷3=279, 2", $true, $false, $false, $false, $false, $true, 1, $false, "553÷8=69, 1", 2) | Out-Null
$d.Content.Find.Execute("922÷7=131, 5", $true, $false, $false, $false, $false, $true, 1, $false, "510÷3=170, 0", 2) | Out-Null
$d.Content.Find.Execute("791÷5=158, 1", $true, $false, $false, $false, $false, $true, 1, $false, "777÷5=155, 2", 2) | Out-Null
$d.Content.Find.Execute("143÷5=28, 3", $true, $false, $false, $false, $false, $true, 1, $false, "805÷9=89, 4", 2) | Out-Null
$d.Content.Find.Execute("335÷8=41, 7", $true, $false, $false, $false, $false, $true, 1, $false, "717÷6=119, 3", 2) | Out-Null
$d.Content.Find.Execute("726÷2=363, 0", $true, $false, $false, $false, $false, $true, 1, $false, "828÷8=103, 4", 2) | Out-Null
